$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-407). The automated export bumped this "last changed" date by
# one day (2023-09-09 -> 2023-09-10, serial 45178 -> 45179) for every row.
for ($r = 2; $r -le 407; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2() + 1
}
